$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1292686666666667
$ws.Range("H2").Value = 0.387806
$ws.Range("M2").Value = 2.815739333333333
$ws.Range("N2").Value = 8.447217999999999
$ws.Range("O2").Value = 0.07700398964630729
$ws.Range("P2").Value = 0.07700398964630729
$ws.Range("Q2").Value = 0.3639868693008889
$ws.Range("R2").Value = 3.275881823708
$ws.Range("S2").Value = 0.07700398964630729
$ws.Range("T2").Value = 0.07700398964630729
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1292686666666667
$ws.Range("H3").Value = 0.387806
$ws.Range("O3").Value = 0.1324338085883186
$ws.Range("P3").Value = 0.1324338085883186
$ws.Range("Q3").Value = 0.6259957126775556
$ws.Range("R3").Value = 5.633961414098
$ws.Range("S3").Value = 0.1324338085883186
$ws.Range("T3").Value = 0.1324338085883186
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1292686666666667
$ws.Range("H4").Value = 0.387806
$ws.Range("M4").Value = 5.537790999999999
$ws.Range("N4").Value = 16.613373
$ws.Range("O4").Value = 0.1514458372546134
$ws.Range("P4").Value = 0.1514458372546134
$ws.Range("Q4").Value = 0.7158628588486666
$ws.Range("R4").Value = 6.442765729637999
$ws.Range("S4").Value = 0.1514458372546134
$ws.Range("T4").Value = 0.1514458372546134
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1292686666666667
$ws.Range("H5").Value = 0.387806
$ws.Range("M5").Value = 1.188595666666667
$ws.Range("N5").Value = 3.565787
$ws.Range("O5").Value = 0.03250535563648733
$ws.Range("P5").Value = 0.03250535563648733
$ws.Range("Q5").Value = 0.1536481770357778
$ws.Range("R5").Value = 1.382833593322
$ws.Range("S5").Value = 0.03250535563648733
$ws.Range("T5").Value = 0.03250535563648733
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1292686666666667
$ws.Range("H6").Value = 0.387806
$ws.Range("M6").Value = 18.85109966666667
$ws.Range("N6").Value = 56.553299
$ws.Range("O6").Value = 0.5155341854158992
$ws.Range("P6").Value = 0.5155341854158992
$ws.Range("Q6").Value = 2.436856519110445
$ws.Range("R6").Value = 21.931708671994
$ws.Range("S6").Value = 0.5155341854158992
$ws.Range("T6").Value = 0.5155341854158992
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1292686666666667
$ws.Range("H7").Value = 0.387806
$ws.Range("M7").Value = 3.330328666666666
$ws.Range("N7").Value = 9.990985999999999
$ws.Range("O7").Value = 0.09107682345837424
$ws.Range("P7").Value = 0.09107682345837424
$ws.Range("Q7").Value = 0.4305071463017778
$ws.Range("R7").Value = 3.874564316716
$ws.Range("S7").Value = 0.09107682345837424
$ws.Range("T7").Value = 0.09107682345837424
